$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A1"
$ws.Range("C2").Value = "M1"
$ws.Range("D2").Value = "A1"
$ws.Range("E2").Value = "A1"
$ws.Range("F2").Value = "DO"
$ws.Range("G2").Value = "M3"
$ws.Range("H2").Value = "M1"
$ws.Range("I2").Value = "DO"
$ws.Range("J2").Value = "M1"
$ws.Range("K2").Value = "M3"
$ws.Range("L2").Value = "M1"
$ws.Range("M2").Value = "M1"
$ws.Range("N2").Value = "M1"
$ws.Range("O2").Value = "M3"
$ws.Range("P2").Value = "A1"
$ws.Range("Q2").Value = "A1"
$ws.Range("R2").Value = "M3"
$ws.Range("S2").Value = "A1"
$ws.Range("T2").Value = "M1"
$ws.Range("U2").Value = "M1"
$ws.Range("V2").Value = "DO"
$ws.Range("W2").Value = "A1"
$ws.Range("X2").Value = "A1"
$ws.Range("Y2").Value = "M1"
$ws.Range("Z2").Value = "DO"
$ws.Range("AA2").Value = "M1"
$ws.Range("AB2").Value = "M3"
$ws.Range("AC2").Value = "M1"
$ws.Range("B3").Value = "M2"
$ws.Range("C3").Value = "A2"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = "DO"
$ws.Range("F3").Value = "M2"
$ws.Range("G3").Value = "M1"
$ws.Range("H3").Value = "A1"
$ws.Range("I3").Value = "A2"
$ws.Range("J3").Value = "M2"
$ws.Range("K3").Value = "A1"
$ws.Range("L3").Value = "A2"
$ws.Range("M3").Value = "DO"
$ws.Range("N3").Value = "M1"
$ws.Range("O3").Value = "A2"
$ws.Range("P3").Value = "DO"
$ws.Range("Q3").Value = "PH"
$ws.Range("R3").Value = "PH"
$ws.Range("S3").Value = "M2"
$ws.Range("T3").Value = "A2"
$ws.Range("U3").Value = "M2"
$ws.Range("V3").Value = "M2"
$ws.Range("W3").Value = "DO"
$ws.Range("X3").Value = "PH"
$ws.Range("Y3").Value = "PH"
$ws.Range("Z3").Value = "A2"
$ws.Range("AA3").Value = "M2"
$ws.Range("AB3").Value = "M2"
$ws.Range("AC3").Value = "M2"
$ws.Range("B4").Value = "DO"
$ws.Range("C4").Value = "M1"
$ws.Range("D4").Value = "M3"
$ws.Range("E4").Value = "M1"
$ws.Range("F4").Value = "A1"
$ws.Range("G4").Value = "A1"
$ws.Range("H4").Value = "M1"
$ws.Range("I4").Value = "DO"
$ws.Range("J4").Value = "A1"
$ws.Range("K4").Value = "M1"
$ws.Range("L4").Value = "M3"
$ws.Range("M4").Value = "A1"
$ws.Range("N4").Value = "A1"
$ws.Range("O4").Value = "M3"
$ws.Range("P4").Value = "M1"
$ws.Range("Q4").Value = "PH"
$ws.Range("R4").Value = "PH"
$ws.Range("S4").Value = "DO"
$ws.Range("T4").Value = "M1"
$ws.Range("U4").Value = "A1"
$ws.Range("V4").Value = "M3"
$ws.Range("W4").Value = "M1"
$ws.Range("X4").Value = "PH"
$ws.Range("Y4").Value = "PH"
$ws.Range("Z4").Value = "DO"
$ws.Range("AA4").Value = "M3"
$ws.Range("AB4").Value = "A1"
$ws.Range("AC4").Value = "M1"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = "A2"
$ws.Range("D5").Value = "A2"
$ws.Range("E5").Value = "M1"
$ws.Range("F5").Value = "M1"
$ws.Range("G5").Value = "M2"
$ws.Range("H5").Value = "DO"
$ws.Range("I5").Value = "M2"
$ws.Range("J5").Value = "M2"
$ws.Range("K5").Value = "A2"
$ws.Range("L5").Value = "A1"
$ws.Range("M5").Value = "DO"
$ws.Range("N5").Value = "M1"
$ws.Range("O5").Value = "M2"
$ws.Range("P5").Value = "M2"
$ws.Range("Q5").Value = "PH"
$ws.Range("R5").Value = "PH"
$ws.Range("S5").Value = "M2"
$ws.Range("T5").Value = "DO"
$ws.Range("U5").Value = "M2"
$ws.Range("V5").Value = "A2"
$ws.Range("W5").Value = "M2"
$ws.Range("X5").Value = "PH"
$ws.Range("Y5").Value = "PH"
$ws.Range("Z5").Value = "M2"
$ws.Range("AA5").Value = "DO"
$ws.Range("AB5").Value = "M2"
$ws.Range("AC5").Value = "A2"
$ws.Range("B6").Value = "DO"
$ws.Range("C6").Value = "M2"
$ws.Range("D6").Value = "A2"
$ws.Range("E6").Value = "M2"
$ws.Range("F6").Value = "A1"
$ws.Range("G6").Value = "M2"
$ws.Range("H6").Value = "M1"
$ws.Range("I6").Value = "M2"
$ws.Range("J6").Value = "M2"
$ws.Range("K6").Value = "A2"
$ws.Range("L6").Value = "A1"
$ws.Range("M6").Value = "M1"
$ws.Range("N6").Value = "M2"
$ws.Range("O6").Value = "DO"
$ws.Range("P6").Value = "DO"
$ws.Range("Q6").Value = "M2"
$ws.Range("R6").Value = "A2"
$ws.Range("S6").Value = "M2"
$ws.Range("T6").Value = "M1"
$ws.Range("U6").Value = "M1"
$ws.Range("V6").Value = "A2"
$ws.Range("W6").Value = "DO"
$ws.Range("X6").Value = "M2"
$ws.Range("Y6").Value = "M2"
$ws.Range("Z6").Value = "M1"
$ws.Range("AA6").Value = "A1"
$ws.Range("AB6").Value = "M2"
$ws.Range("AC6").Value = "A2"
$ws.Range("B7").Value = "A1"
$ws.Range("C7").Value = "A1"
$ws.Range("D7").Value = "M3"
$ws.Range("E7").Value = "DO"
$ws.Range("F7").Value = "M1"
$ws.Range("G7").Value = "A1"
$ws.Range("H7").Value = "A1"
$ws.Range("I7").Value = "A1"
$ws.Range("J7").Value = "A1"
$ws.Range("K7").Value = "DO"
$ws.Range("L7").Value = "M3"
$ws.Range("M7").Value = "A1"
$ws.Range("N7").Value = "A1"
$ws.Range("O7").Value = "A1"
$ws.Range("P7").Value = "A1"
$ws.Range("Q7").Value = "A1"
$ws.Range("R7").Value = "M3"
$ws.Range("S7").Value = "A1"
$ws.Range("T7").Value = "A1"
$ws.Range("U7").Value = "A1"
$ws.Range("V7").Value = "DO"
$ws.Range("W7").Value = "M1"
$ws.Range("X7").Value = "M3"
$ws.Range("Y7").Value = "A1"
$ws.Range("Z7").Value = "A1"
$ws.Range("AA7").Value = "DO"
$ws.Range("AB7").Value = "A1"
$ws.Range("AC7").Value = "A1"
$ws.Range("B8").Value = "A2"
$ws.Range("C8").Value = "A2"
$ws.Range("D8").Value = "DO"
$ws.Range("E8").Value = "A2"
$ws.Range("F8").Value = "A1"
$ws.Range("G8").Value = "M1"
$ws.Range("H8").Value = "A2"
$ws.Range("I8").Value = "A1"
$ws.Range("J8").Value = "A1"
$ws.Range("K8").Value = "M3"
$ws.Range("L8").Value = "A1"
$ws.Range("M8").Value = "A1"
$ws.Range("N8").Value = "M1"
$ws.Range("O8").Value = "DO"
$ws.Range("P8").Value = "A2"
$ws.Range("Q8").Value = "M2"
$ws.Range("R8").Value = "A2"
$ws.Range("S8").Value = "A2"
$ws.Range("T8").Value = "DO"
$ws.Range("U8").Value = "M1"
$ws.Range("V8").Value = "M1"
$ws.Range("W8").Value = "A2"
$ws.Range("X8").Value = "A2"
$ws.Range("Y8").Value = "M2"
$ws.Range("Z8").Value = "A2"
$ws.Range("AA8").Value = "A1"
$ws.Range("AB8").Value = "M1"
$ws.Range("AC8").Value = "DO"
$ws.Range("B9").Value = "M2"
$ws.Range("C9").Value = "M2"
$ws.Range("D9").Value = "A1"
$ws.Range("E9").Value = "A2"
$ws.Range("F9").Value = "M1"
$ws.Range("G9").Value = "M2"
$ws.Range("H9").Value = "DO"
$ws.Range("I9").Value = "DO"
$ws.Range("J9").Value = "M2"
$ws.Range("K9").Value = "M2"
$ws.Range("L9").Value = "M2"
$ws.Range("M9").Value = "M2"
$ws.Range("N9").Value = "A1"
$ws.Range("O9").Value = "M1"
$ws.Range("P9").Value = "M2"
$ws.Range("Q9").Value = "A2"
$ws.Range("R9").Value = "M2"
$ws.Range("S9").Value = "M2"
$ws.Range("T9").Value = "A1"
$ws.Range("U9").Value = "A1"
$ws.Range("V9").Value = "DO"
$ws.Range("W9").Value = "A2"
$ws.Range("X9").Value = "M2"
$ws.Range("Y9").Value = "A1"
$ws.Range("Z9").Value = "M2"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "A1"
$ws.Range("AC9").Value = "DO"
$ws.Range("B10").Value = "DO"
$ws.Range("C10").Value = "M2"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = "M2"
$ws.Range("F10").Value = "A2"
$ws.Range("G10").Value = "A1"
$ws.Range("H10").Value = "M1"
$ws.Range("I10").Value = "DO"
$ws.Range("J10").Value = "M2"
$ws.Range("K10").Value = "A2"
$ws.Range("L10").Value = "A1"
$ws.Range("M10").Value = "M2"
$ws.Range("N10").Value = "M2"
$ws.Range("O10").Value = "A1"
$ws.Range("P10").Value = "M2"
$ws.Range("Q10").Value = "PH"
$ws.Range("R10").Value = "PH"
$ws.Range("S10").Value = "DO"
$ws.Range("T10").Value = "M2"
$ws.Range("U10").Value = "M2"
$ws.Range("V10").Value = "A2"
$ws.Range("W10").Value = "M2"
$ws.Range("X10").Value = "PH"
$ws.Range("Y10").Value = "PH"
$ws.Range("Z10").Value = "M2"
$ws.Range("AA10").Value = "A2"
$ws.Range("AB10").Value = "M2"
$ws.Range("AC10").Value = "DO"
